# Applies the changes described by the commit "Add files via upload":
#   1. The "Duración" label textbox (Google Shape;89;p1) is repositioned /
#      resized slightly and its paragraph is re-centred.
#   2. The old dashed-border "QR AQUÍ" placeholder (Google Shape;92;p1) is removed.
#   3. A new plain textbox named "CuadroTexto 2" with centred "QR AQUÍ" text
#      is added in roughly the same spot.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Reposition / resize the "Dur..." label and centre its text ------
# (Only Left/Width actually change in the target; Top/Height are left
#  untouched so their EMU values survive the point round-trip exactly.)
$durShape = $s.Shapes.Item("Google Shape;89;p1")
$durShape.Left = 526.7143307086615
$durShape.Width = 51.85716535433071
$durShape.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# --- 2. Remove the old dashed rectangle "QR AQUÍ" placeholder shape -----
$oldQr = $s.Shapes.Item("Google Shape;92;p1")
$oldQr.Delete()

# --- 3. Add the new "CuadroTexto 2" textbox with centred QR AQUÍ text ---
$qr = $s.Shapes.AddTextbox(1, 423.3070866141732, 413.1427559055118, 113.3858267716535, 113.3858267716535)
$qr.Name = "CuadroTexto 2"
$qr.Fill.Visible = 0
$qr.TextFrame.WordWrap = 0
$qr.TextFrame.AutoSize = 0
$qr.TextFrame.VerticalAnchor = 3
$qr.TextFrame.TextRange.Text = "QR AQUÍ"
$qr.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$qr.TextFrame.TextRange.LanguageID = "es-MX"
